{"js": "// Update the date heading and every two-digit-by-two-digit multiplication\n// answer in the table with the new values from the latest run.\n\n// 1) Update the date/weekday heading paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstPara = paragraphs.items[0];\nfirstPara.load(\"text\");\nawait context.sync();\nif (firstPara.text.trim() === \"2024-11-24 Sunday\") {\n  // Replace the run's text in place (via its range) so the existing\n  // character formatting (Arial, size 30) on the run is preserved.\n  const headingRange = firstPara.getRange();\n  headingRange.insertText(\"2024-11-25 Monday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the multiplication table's cell contents in place, preserving\n// each cell's existing formatting.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Old answer -> new answer, keyed by the exact text that was in the cell.\nconst replacements = {\n  \"49\u00d715=735\": \"44\u00d711=484\",\n  \"15\u00d740=600\": \"92\u00d785=7820\",\n  \"83\u00d716=1328\": \"81\u00d733=2673\",\n  \"37\u00d760=2220\": \"14\u00d735=490\",\n  \"48\u00d758=2784\": \"70\u00d780=5600\",\n  \"46\u00d739=1794\": \"80\u00d740=3200\",\n  \"85\u00d792=7820\": \"78\u00d711=858\",\n  \"28\u00d737=1036\": \"17\u00d726=442\",\n  \"70\u00d748=3360\": \"45\u00d779=3555\",\n  \"67\u00d734=2278\": \"44\u00d791=4004\",\n  \"53\u00d768=3604\": \"34\u00d766=2244\",\n  \"83\u00d757=4731\": \"90\u00d719=1710\",\n  \"12\u00d752=624\": \"56\u00d737=2072\",\n  \"41\u00d777=3157\": \"40\u00d771=2840\",\n  \"89\u00d779=7031\": \"69\u00d765=4485\",\n  \"19\u00d737=703\": \"51\u00d745=2295\",\n  \"59\u00d755=3245\": \"27\u00d773=1971\",\n  \"13\u00d740=520\": \"70\u00d744=3080\",\n  \"45\u00d781=3645\": \"62\u00d777=4774\",\n  \"64\u00d774=4736\": \"42\u00d746=1932\",\n  \"52\u00d722=1144\": \"55\u00d792=5060\",\n  \"54\u00d731=1674\": \"24\u00d740=960\",\n  \"55\u00d743=2365\": \"28\u00d730=840\",\n  \"53\u00d758=3074\": \"97\u00d787=8439\",\n  \"85\u00d783=7055\": \"64\u00d718=1152\",\n};\n\nconst values = table.values;\nfor (let r = 0; r < values.length; r++) {\n  for (let c = 0; c < values[r].length; c++) {\n    const oldText = values[r][c];\n    if (Object.prototype.hasOwnProperty.call(replacements, oldText)) {\n      const cell = table.getCell(r, c);\n      cell.value = replacements[oldText];\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date heading and every two-digit-by-two-digit multiplication\n# answer in the table with the new values from the latest run.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday heading paragraph (first paragraph in the body).\n$headingPara = $d.Paragraphs.Item(1)\nif ($headingPara.Range.Text.Trim() -eq \"2024-11-24 Sunday\") {\n    $headingPara.Range.Text = \"2024-11-25 Monday\"\n}\n\n# 2) Update the multiplication table's cell contents in place, preserving\n# each cell's existing formatting.\n$table = $d.Tables.Item(1)\n\n# Old answer -> new answer, keyed by the exact text that was in the cell.\n$replacements = @{\n    \"49\u00d715=735\"  = \"44\u00d711=484\"\n    \"15\u00d740=600\"  = \"92\u00d785=7820\"\n    \"83\u00d716=1328\" = \"81\u00d733=2673\"\n    \"37\u00d760=2220\" = \"14\u00d735=490\"\n    \"48\u00d758=2784\" = \"70\u00d780=5600\"\n    \"46\u00d739=1794\" = \"80\u00d740=3200\"\n    \"85\u00d792=7820\" = \"78\u00d711=858\"\n    \"28\u00d737=1036\" = \"17\u00d726=442\"\n    \"70\u00d748=3360\" = \"45\u00d779=3555\"\n    \"67\u00d734=2278\" = \"44\u00d791=4004\"\n    \"53\u00d768=3604\" = \"34\u00d766=2244\"\n    \"83\u00d757=4731\" = \"90\u00d719=1710\"\n    \"12\u00d752=624\"  = \"56\u00d737=2072\"\n    \"41\u00d777=3157\" = \"40\u00d771=2840\"\n    \"89\u00d779=7031\" = \"69\u00d765=4485\"\n    \"19\u00d737=703\"  = \"51\u00d745=2295\"\n    \"59\u00d755=3245\" = \"27\u00d773=1971\"\n    \"13\u00d740=520\"  = \"70\u00d744=3080\"\n    \"45\u00d781=3645\" = \"62\u00d777=4774\"\n    \"64\u00d774=4736\" = \"42\u00d746=1932\"\n    \"52\u00d722=1144\" = \"55\u00d792=5060\"\n    \"54\u00d731=1674\" = \"24\u00d740=960\"\n    \"55\u00d743=2365\" = \"28\u00d730=840\"\n    \"53\u00d758=3074\" = \"97\u00d787=8439\"\n    \"85\u00d783=7055\" = \"64\u00d718=1152\"\n}\n\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        $cell = $table.Cell($r, $c)\n        # A table cell's Range.Text includes a trailing end-of-cell marker\n        # (CR + BEL, i.e. chars 13/7) that isn't part of the visible text.\n        $rawText = $cell.Range.Text\n        $cellText = $rawText.Substring(0, $rawText.Length - 2)\n        if ($replacements.ContainsKey($cellText)) {\n            $cell.Range.Text = $replacements[$cellText]\n        }\n    }\n}\n"}
